# Apply "repull data, push all data, mean calculation" updates to the dSF
# column (column F) for a handful of rows in the pitcher log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of worksheet row -> new dSF (column F) value
$updates = @{
    5  = 3
    13 = 1
    16 = 1
    22 = 1
    23 = 1
    26 = 3
    30 = -1
    34 = 1
    41 = 0
    48 = 1
    52 = 3
    57 = 0
    58 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
